$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordered data (language, value) - sorted descending by value, with
# "Swedish" and "Uzbek" removed from the original list.
$data = @(
    @("Chinese", 20.95281307444563),
    @("English", 20.08641337179922),
    @("Spanish", 5.91422511873983),
    @("Arabic", 4.236767452351714),
    @("German", 3.753680209417977),
    @("Malay-Indonesian", 3.453750343523509),
    @("Japanese", 3.36126385077831),
    @("Russian", 2.747580672356291),
    @("Portuguese", 2.647227169817934),
    @("French", 2.366244602578518),
    @("Turkish", 2.053916697825131),
    @("Italian", 1.691257136835023),
    @("Korean", 1.627770076329035),
    @("Dutch", 1.10876509614138),
    @("Bengali", 0.9844339970642147),
    @("Polish", 0.9811041279136866),
    @("Vietnamese", 0.9563494208700318),
    @("Urdu", 0.9561181211916385),
    @("Persian", 0.9317202362128222),
    @("Thai", 0.9182268147306173)
)

# Clear out the old data rows (2 through 23) entirely, keeping the header row intact.
$oldLastRow = 23
$ws.Range("A2:B$oldLastRow").Clear()

# Write the new data back starting at row 2, in the new sorted order.
$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

$newLastRow = $row - 1

# Re-apply the language-column formatting that the header cell (A1) carries
# (bold font, thin box border, centered/top aligned) to every language cell,
# matching the original style used throughout column A.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A$newLastRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
